$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-07-17 Monday" "2023-07-18 Tuesday"
Replace-Text "30-26=4" "31+37=68"
Replace-Text "52+35=87" "95-45=50"
Replace-Text "80-46=34" "44-15=29"
Replace-Text "73-49=24" "77-44=33"
Replace-Text "83+7=90" "59-14=45"
Replace-Text "74-28=46" "2+35=37"
Replace-Text "25-23=2" "6+88=94"
Replace-Text "14+12=26" "1+68=69"
Replace-Text "40+10=50" "89-49=40"
Replace-Text "57-42=15" "53+18=71"
Replace-Text "28+44=72" "97-53=44"
Replace-Text "10+74=84" "69-1=68"
Replace-Text "95-68=27" "67+4=71"
Replace-Text "72-25=47" "93-89=4"
Replace-Text "24+21=45" "37+33=70"
Replace-Text "36+13=49" "79-75=4"
Replace-Text "26-12=14" "50-47=3"
Replace-Text "19+72=91" "33+12=45"
Replace-Text "59-12=47" "51-3=48"
Replace-Text "80-36=44" "85-57=28"
Replace-Text "83-57=26" "70-9=61"
Replace-Text "17+1=18" "9+24=33"
Replace-Text "41+5=46" "52-41=11"
Replace-Text "63-54=9" "29+4=33"
Replace-Text "80-41=39" "52-7=45"
Replace-Text "21+62=83" "40+48=88"
Replace-Text "71+24=95" "28+68=96"
Replace-Text "68-39=29" "72+16=88"
Replace-Text "32+61=93" "36+22=58"
Replace-Text "44+7=51" "30+67=97"
Replace-Text "80-45=35" "14+39=53"
Replace-Text "42+43=85" "71-64=7"
Replace-Text "45+0=45" "25-0=25"
Replace-Text "76-58=18" "37+8=45"
Replace-Text "3+83=86" "2+3=5"
Replace-Text "31+45=76" "40+23=63"
Replace-Text "93-3=90" "77+14=91"
Replace-Text "72-5=67" "34-22=12"
Replace-Text "86-15=71" "75-68=7"
Replace-Text "51+16=67" "96-77=19"
Replace-Text "54+0=54" "35-17=18"
Replace-Text "55-48=7" "5+77=82"
Replace-Text "85-60=25" "63-39=24"
Replace-Text "65+5=70" "20+35=55"
Replace-Text "11+51=62" "22+6=28"
Replace-Text "86-1=85" "98-80=18"
Replace-Text "66-63=3" "82+1=83"
Replace-Text "78-27=51" "45+1=46"
Replace-Text "81-46=35" "17+68=85"
Replace-Text "97-87=10" "23+47=70"
Replace-Text "52+32=84" "29+64=93"
Replace-Text "1+75=76" "96-57=39"
Replace-Text "20+16=36" "8+47=55"
Replace-Text "28+36=64" "39-2=37"
Replace-Text "95-42=53" "73-57=16"
Replace-Text "39+11=50" "36+9=45"
Replace-Text "71-53=18" "7+22=29"
Replace-Text "6+49=55" "42+55=97"
Replace-Text "92-83=9" "18+13=31"
Replace-Text "33-4=29" "77+3=80"
Replace-Text "68-5=63" "23-22=1"
Replace-Text "97-71=26" "99-39=60"
Replace-Text "82-21=61" "82-28=54"
Replace-Text "66+5=71" "9+4=13"
Replace-Text "86-25=61" "46-2=44"
Replace-Text "30+47=77" "36+1=37"
Replace-Text "1+89=90" "97-95=2"
Replace-Text "35-0=35" "78+17=95"
Replace-Text "80-12=68" "70-19=51"
Replace-Text "81-48=33" "74+20=94"
Replace-Text "64+16=80" "56-47=9"
Replace-Text "19+58=77" "89-62=27"
Replace-Text "84-53=31" "49-42=7"
Replace-Text "45-31=14" "26-16=10"
Replace-Text "86-7=79" "25+24=49"
Replace-Text "27+67=94" "59+36=95"
Replace-Text "33+4=37" "69+17=86"
Replace-Text "93-79=14" "15+71=86"
Replace-Text "35-16=19" "0+89=89"
Replace-Text "77+6=83" "79-54=25"
Replace-Text "71-54=17" "85-51=34"
Replace-Text "90-13=77" "2+81=83"
Replace-Text "7+70=77" "68+4=72"
Replace-Text "32+38=70" "99-93=6"
Replace-Text "80+4=84" "45+9=54"
Replace-Text "37+17=54" "67+32=99"
Replace-Text "60-53=7" "46+26=72"
Replace-Text "3+10=13" "57+33=90"
Replace-Text "34+52=86" "15+28=43"
Replace-Text "83-26=57" "71-7=64"
Replace-Text "17-4=13" "69-51=18"
Replace-Text "5+91=96" "90-19=71"
Replace-Text "99-23=76" "6+86=92"
Replace-Text "81-45=36" "53-7=46"
Replace-Text "13+42=55" "93-23=70"
Replace-Text "74-11=63" "62+6=68"
Replace-Text "40+11=51" "90-70=20"
Replace-Text "51+5=56" "42+41=83"
Replace-Text "29+7=36" "6+19=25"
Replace-Text "51+17=68" "54-6=48"
